# upbit auto trade v.2.2
# Adds the newest trade log rows to Sheet1, the matching price snapshot row
# to Sheet2, and refreshes the timestamp recorded for the last "거미줄매수"
# batch on Sheet3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1 : trade history - two new executions on 2022-02-15
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# Seed rows 51:52 from row 50 so the new rows inherit the same (default)
# cell styling as the rest of the table instead of the raw column style.
$ws1.Range("A50:F50").Copy($ws1.Range("A51:F52"))

# Row 51 - KRW-XRP sell
$ws1.Cells.Item(51, 1).Value = "'2022-02-15"
$ws1.Cells.Item(51, 2).Value = "10:21:59"
$ws1.Cells.Item(51, 3).Value = "KRW-XRP"
$ws1.Cells.Item(51, 4).Value = "매도"
$ws1.Cells.Item(51, 5).Value = "'11.67606225"
$ws1.Cells.Item(51, 6).Value = "9cd6b065-59c6-4079-9d01-6552963a6f63"

# Row 52 - KRW-BTC buy
$ws1.Cells.Item(52, 1).Value = "'2022-02-15"
$ws1.Cells.Item(52, 2).Value = "23:52:01"
$ws1.Cells.Item(52, 3).Value = "KRW-BTC"
$ws1.Cells.Item(52, 4).Value = "매수"
$ws1.Cells.Item(52, 5).Value = "'10604.0"
$ws1.Cells.Item(52, 6).Value = "62f5a1a6-e346-47a6-bf80-a5850eec89c5"

# ---------------------------------------------------------------------
# Sheet2 : ticker / target price snapshots - new KRW-BTC row
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

# Seed row 93 from row 92 so it keeps the plain (default) numeric styling
# plus the shared datetime style used in column D.
$ws2.Range("A92:E92").Copy($ws2.Range("A93:E93"))

$ws2.Cells.Item(93, 1).Value = "KRW-BTC"
$ws2.Cells.Item(93, 2).Value = 52234500
$ws2.Cells.Item(93, 3).Value = 51947600
$ws2.Cells.Item(93, 4).Value = 44607.99442908435
$ws2.Cells.Item(93, 5).Value = 10604

# ---------------------------------------------------------------------
# Sheet3 : 거미줄매수 order log - refresh the shared timestamp used by the
# last batch (rows 35-42) to the recalculated value.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")

for ($r = 35; $r -le 42; $r++) {
    $ws3.Cells.Item($r, 2).Value = 44606.90676508102
}
